$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the coin/link/price/volume text columns remain plain text so that
# values such as "1.004" or "0.09500" are not auto-converted to numbers
# (which would silently drop significant trailing zeros).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.458.88'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.902.07'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '325.42'
$ws.Range("E5").Value = '  -2.41%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = '0.4791'
$ws.Range("E7").Value = '  +2.60%  '
$ws.Range("D8").Value = '0.4062'
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("D9").Value = '0.08072'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").Value = '23.39'
$ws.Range("E11").Value = '  +4.52%  '
$ws.Range("D12").Value = '1.922.45'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '5.955'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '7.076'
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = '90.05'
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("D17").Value = '0.06689'
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '17.63'
$ws.Range("E19").Value = '  -0.90%  '
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '29.481.95'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").Value = '5.539'
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = '11.81'
$ws.Range("E23").Value = '  +2.52%  '
$ws.Range("D24").Value = '2.167'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '2.142.39'
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("D26").Value = '154.36'
$ws.Range("E26").Value = '  -0.54%  '
$ws.Range("D27").Value = '19.83'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '6.099'
$ws.Range("E28").Value = '  +5.91%  '
$ws.Range("D29").Value = '2.093'
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").Value = '118.43'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").Value = '1.031'
$ws.Range("E31").Value = '  -3.65%  '
$ws.Range("D32").Value = '0.09500'
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").Value = '5.470'
$ws.Range("E33").Value = '  +1.25%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.545'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.390'
$ws.Range("E35").Value = '  -2.51%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.06074'
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02250'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("D38").Value = '1.173'
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("D39").Value = '0.5891'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").Value = '7.917'
$ws.Range("E40").Value = '  -5.94%  '
$ws.Range("D41").Value = '0.1844'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").Value = '10.18'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = '1.283'
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = '0.07767'
$ws.Range("E44").Value = '  +3.40%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '2.388'
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("D46").Value = '12.27'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").Value = '0.5529'
$ws.Range("D48").Value = '1.923'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").Value = '113.90'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("D50").Value = '72.35'
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").Value = '0.2938'
$ws.Range("E51").Value = '  -1.34%  '
